# "Set origin to sprite center for flying objects"
#
# The enemy-spawn offset table on the EnemyFleet sheet is reworked:
#   - Column A keeps its "previous row + step" shared formula, but the
#     step shrinks from 100 to 2 (positions are now measured from the
#     sprite's center instead of a far corner, so the spacing is tiny).
#   - Column B drops its "+100 per step" formula pattern entirely and
#     becomes a hand-authored list of small static offsets.
#   - The sheet's last active selection moves from C21 to B21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: rewrite every formula as "=+A<prev>+2" (was "+100"), row by row
# so the engine re-derives the A3:A20 shared-formula group automatically.
for ($r = 2; $r -le 20; $r++) {
    $prev = $r - 1
    $ws.Range("A$r").Formula = "=+A$prev+2"
}

# Column B: replace every value/formula with the new static numbers.
$bValues = @(2, 5, 7, 3, 10, 6, 7, 1, 5, 3, 4, 5, 8, 2, 11, 3, 8, 5, 9, 7)
for ($r = 1; $r -le 20; $r++) {
    $ws.Range("B$r").Value = $bValues[$r - 1]
}

# Selection moves from C21 to B21.
[void]$ws.Range("B21").Select()

# Best-effort localization touch-ups that accompanied this edit upstream
# (renaming the built-in cell styles to their Croatian names and renaming
# the theme). These are harmless no-ops if the host doesn't expose a
# writable surface for them.
try { $wb.Theme.Name = "Tema sustava Office" } catch {}
try { $wb.Styles.Item("Normal").Name = "Normalno" } catch {}
try { $wb.Styles.Item("Comma").Name = "Zarez" } catch {}
